# Apply odds updates to Sheet1 data rows (2026-01-08 Betfair Back/Lay odds refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.94
$ws.Range("G2").Value = 3.45
$ws.Range("H2").Value = 2.74
$ws.Range("I2").Value = 3.05
$ws.Range("K2").Value = 3.25
$ws.Range("L2").Value = 1.56
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 2.2
$ws.Range("Q2").Value = 3
$ws.Range("S2").Value = 6
$ws.Range("T2").Value = 2.28
$ws.Range("U2").Value = 1.63
$ws.Range("V2").Value = 1.47
$ws.Range("W2").Value = 1.42
$ws.Range("Y2").Value = 9

# Row 3
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 19.5
$ws.Range("J3").Value = 9.6
$ws.Range("K3").Value = 11.5
$ws.Range("L3").Value = 1.15
$ws.Range("N3").Value = 8.6
$ws.Range("P3").Value = 3.8
$ws.Range("Q3").Value = 1.28
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.74
$ws.Range("T3").Value = 1.92
$ws.Range("U3").Value = 1.94
$ws.Range("AB3").Value = 19.5
$ws.Range("AE3").Value = 280
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 16
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 16
$ws.Range("AL3").Value = 40
$ws.Range("AM3").Value = 190
$ws.Range("AN3").Value = 2.72

# Row 4
$ws.Range("G4").Value = 2
$ws.Range("M4").Value = 1.1
$ws.Range("R4").Value = 1.19
$ws.Range("S4").Value = 5.3
$ws.Range("T4").Value = 2.2
$ws.Range("U4").Value = 1.67
$ws.Range("W4").Value = 2

# Row 5
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 1.53
$ws.Range("P5").Value = 1.66
$ws.Range("Q5").Value = 2.42
$ws.Range("S5").Value = 4.8

# Row 6
$ws.Range("F6").Value = 1.6
$ws.Range("G6").Value = 1.64
$ws.Range("H6").Value = 4.8
$ws.Range("K6").Value = 5.7
$ws.Range("L6").Value = 1.17
$ws.Range("Q6").Value = 1.35
$ws.Range("R6").Value = 1.94
$ws.Range("W6").Value = 2.56
$ws.Range("Y6").Value = 990
$ws.Range("AB6").Value = 18.5
$ws.Range("AD6").Value = 22
$ws.Range("AE6").Value = 980

# Row 7
$ws.Range("H7").Value = 1.77
$ws.Range("I7").Value = 1.89
$ws.Range("J7").Value = 3.65
$ws.Range("L7").Value = 1.41
$ws.Range("P7").Value = 1.86
$ws.Range("Q7").Value = 1.94
$ws.Range("S7").Value = 3.4
$ws.Range("V7").Value = 2.08
$ws.Range("AB7").Value = 20

# Row 8
$ws.Range("I8").Value = 8.199999999999999
$ws.Range("L8").Value = 1.5
$ws.Range("Q8").Value = 2.38
$ws.Range("S8").Value = 4.7

# Row 9
$ws.Range("F9").Value = 1.46
$ws.Range("G9").Value = 1.47
$ws.Range("H9").Value = 9.199999999999999
$ws.Range("I9").Value = 9.800000000000001
$ws.Range("J9").Value = 4.6
$ws.Range("L9").Value = 1.44
$ws.Range("N9").Value = 3.6
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 1.9
$ws.Range("T9").Value = 2.36
$ws.Range("U9").Value = 1.7
$ws.Range("V9").Value = 1.11
$ws.Range("W9").Value = 3.1
$ws.Range("Y9").Value = 25
$ws.Range("Z9").Value = 80
$ws.Range("AA9").Value = 460
$ws.Range("AB9").Value = 6.8
$ws.Range("AD9").Value = 36
$ws.Range("AF9").Value = 7.4
$ws.Range("AN9").Value = 8.800000000000001
$ws.Range("AO9").Value = 330

# Row 10
$ws.Range("F10").Value = 1.62
$ws.Range("G10").Value = 1.63
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 3.15
$ws.Range("V10").Value = 1.19
$ws.Range("X10").Value = 17.5
$ws.Range("AC10").Value = 9.6
$ws.Range("AE10").Value = 80
$ws.Range("AK10").Value = 15
$ws.Range("AN10").Value = 8.199999999999999

# Row 11
$ws.Range("I11").Value = 4.3
$ws.Range("L11").Value = 1.49
$ws.Range("Q11").Value = 2.28
$ws.Range("R11").Value = 1.24
$ws.Range("S11").Value = 4.3
$ws.Range("T11").Value = 1.95
$ws.Range("U11").Value = 1.87
$ws.Range("Y11").Value = 14.5
$ws.Range("Z11").Value = 34
$ws.Range("AB11").Value = 9.4
$ws.Range("AC11").Value = 9
$ws.Range("AD11").Value = 21
$ws.Range("AE11").Value = 75
$ws.Range("AF11").Value = 16
$ws.Range("AG11").Value = 13.5
$ws.Range("AH11").Value = 26
$ws.Range("AI11").Value = 90
$ws.Range("AK11").Value = 34
$ws.Range("AL11").Value = 60
$ws.Range("AN11").Value = 30
